$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.429.63"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -5.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.899.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.71"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.890.43"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -8.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -10.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("E13").Value = "  -5.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.380.13"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.898.96"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.56"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.443.45"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "404.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.668"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.80"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.87"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.65"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0970"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.906"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.41"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -11.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.93"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0616"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.67%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0340"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.616.53"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "359.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "119.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.228"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.83"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.97%  "
